$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price rows arrived; insert them above the existing data
# (which starts at row 51) so every prior record shifts down by two rows.
$ws.Rows.Item(51).Insert()
$ws.Rows.Item(51).Insert()

# Common/constant columns shared by every data row in this sheet.
$commonA = 1
$commonB = "Agrícola del Norte S.A. de Arica"
$commonC = "Arica y Parinacota"
$commonE = 15
$commonF = "Fruta"
$commonG = 100103
$commonH = "Frutos de hueso (carozo)"
$commonI = 100103006
$commonJ = "Nectarín"
$commonT = 18

# Row 51: new "Artic Snow" record
$ws.Cells.Item(51, 1).Value = $commonA
$ws.Cells.Item(51, 2).Value = $commonB
$ws.Cells.Item(51, 3).Value = $commonC
$ws.Cells.Item(51, 4).Value = 44637
$ws.Cells.Item(51, 5).Value = $commonE
$ws.Cells.Item(51, 6).Value = $commonF
$ws.Cells.Item(51, 7).Value = $commonG
$ws.Cells.Item(51, 8).Value = $commonH
$ws.Cells.Item(51, 9).Value = $commonI
$ws.Cells.Item(51, 10).Value = $commonJ
$ws.Cells.Item(51, 11).Value = "Artic Snow"
$ws.Cells.Item(51, 12).Value = "Segunda"
$ws.Cells.Item(51, 13).Value = 300
$ws.Cells.Item(51, 14).Value = 18000
$ws.Cells.Item(51, 15).Value = 20000
$ws.Cells.Item(51, 16).Value = 19000
$ws.Cells.Item(51, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(51, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(51, 19).Value = 1056
$ws.Cells.Item(51, 20).Value = $commonT

# Row 52: new "August Red" record
$ws.Cells.Item(52, 1).Value = $commonA
$ws.Cells.Item(52, 2).Value = $commonB
$ws.Cells.Item(52, 3).Value = $commonC
$ws.Cells.Item(52, 4).Value = 44637
$ws.Cells.Item(52, 5).Value = $commonE
$ws.Cells.Item(52, 6).Value = $commonF
$ws.Cells.Item(52, 7).Value = $commonG
$ws.Cells.Item(52, 8).Value = $commonH
$ws.Cells.Item(52, 9).Value = $commonI
$ws.Cells.Item(52, 10).Value = $commonJ
$ws.Cells.Item(52, 11).Value = "August Red"
$ws.Cells.Item(52, 12).Value = "Segunda"
$ws.Cells.Item(52, 13).Value = 300
$ws.Cells.Item(52, 14).Value = 19000
$ws.Cells.Item(52, 15).Value = 20000
$ws.Cells.Item(52, 16).Value = 19500
$ws.Cells.Item(52, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(52, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(52, 19).Value = 1083
$ws.Cells.Item(52, 20).Value = $commonT
